$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.140.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.145.69"
$ws.Range("D3").Style = "Normal"

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.131.91"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.81%  "

$ws.Range("E10").Value = "  +1.70%  "

$ws.Range("E11").Value = "  +5.52%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.665.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.47%  "

$ws.Range("E16").Value = "  -0.30%  "

$ws.Range("E17").Value = "  +2.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.935.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.135.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "469.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.69%  "

$ws.Range("E21").Value = "  +1.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.731"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +13.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "80.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.38%  "

$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.00%  "

$ws.Range("E29").Value = "  +1.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.46%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.65%  "

$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("E33").Value = "  +4.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0855"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.19%  "

$ws.Range("E36").Value = "  +3.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.91%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "459.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.10%  "

$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.290"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0373"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.888.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.05%  "

$ws.Range("E47").Value = "  -0.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.70%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("E50").Value = "  +0.62%  "

$ws.Range("E51").Value = "  +3.45%  "
